$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URI")

# Row 4 - Inventory
$ws.Range("C4").Value = 125000000.0
$ws.Range("E4").Value = 107000000.0
$ws.Range("F4").Value = 115000000.0
$ws.Range("G4").Value = 120000000.0

# Row 12 - Accounts Payable
$ws.Range("C12").Value = 466000000.0
$ws.Range("D12").Value = 541000000.0
$ws.Range("E12").Value = 316000000.0
$ws.Range("F12").Value = 484000000.0
$ws.Range("G12").Value = 454000000.0

# Row 18 - Long Term Tax Liability (Deferred)
$ws.Range("C18").Value = 1768000000.0
$ws.Range("D18").Value = 1818000000.0
$ws.Range("E18").Value = 1820000000.0
$ws.Range("F18").Value = 1878000000.0
$ws.Range("G18").Value = 1887000000.0
